$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $rowA, $rowB, $firstCol, $lastCol) {
    $rangeA = $ws.Range("$firstCol$rowA`:$lastCol$rowA")
    $rangeB = $ws.Range("$firstCol$rowB`:$lastCol$rowB")

    $valuesA = $rangeA.Value()
    $valuesB = $rangeB.Value()

    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}

# Swap data (columns B through AD) between row 99 and row 100,
# keeping column A (id) untouched.
Swap-Rows $ws 99 100 "B" "AD"

# Swap data (columns B through AD) between row 177 and row 178,
# keeping column A (id) untouched.
Swap-Rows $ws 177 178 "B" "AD"
